$d = $word.ActiveDocument

# The paragraph "<id>p141r_1</id>" is currently split across three runs with
# different character formatting (Courier New / Arial / Courier New). The
# edit merges them into a single run using the first run's formatting
# (Courier New, color 7f6000, sz 18) while keeping the same visible text.
$rng = $d.Content
$found = $rng.Find.Execute("<id>p141r_1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p141r_1</id>", 2)
Write-Host "Replaced:" $found
